# Auto-generated Excel COM-interop script applying diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 788.7273
$ws.Range("I28").Value = 427.9375
$ws.Range("J28").Value = 1750.8334
$ws.Range("K28").Value = 427.9375
$ws.Range("L28").Value = 1750.8334
$ws.Range("M28").Value = 57.0625
$ws.Range("N28").Value = -2720.8334
$ws.Range("H86").Value = 4966.6665
$ws.Range("I86").Value = 3600
$ws.Range("J86").Value = 6333.3335
$ws.Range("K86").Value = 3600
$ws.Range("L86").Value = 6333.3335
$ws.Range("M86").Value = -2477
$ws.Range("N86").Value = -8579.333500000001
$ws.Range("H89").Value = 4966.6665
$ws.Range("I89").Value = 3600
$ws.Range("J89").Value = 6333.3335
$ws.Range("K89").Value = 18000
$ws.Range("L89").Value = 31666.6675
$ws.Range("M89").Value = -12384
$ws.Range("N89").Value = -42898.6675
$ws.Range("H92").Value = 808.94116
$ws.Range("I92").Value = 841.44446
$ws.Range("J92").Value = 683.5714
$ws.Range("K92").Value = 841.44446
$ws.Range("L92").Value = 683.5714
$ws.Range("M92").Value = 406.55554
$ws.Range("N92").Value = -3179.5714
$ws.Range("H121").Value = 942.23254
$ws.Range("J121").Value = 982.9
$ws.Range("L121").Value = 2948.7
$ws.Range("N121").Value = -6442.7
$ws.Range("H132").Value = 5104965
$ws.Range("I132").Value = 2959.5642
$ws.Range("J132").Value = 25002786
$ws.Range("K132").Value = 8878.692599999998
$ws.Range("L132").Value = 75008358
$ws.Range("M132").Value = -6348.692599999998
$ws.Range("N132").Value = -75013418
$ws.Range("H135").Value = 1804.9678
$ws.Range("I135").Value = 1777.6957
$ws.Range("J135").Value = 1883.375
$ws.Range("K135").Value = 15999.2613
$ws.Range("L135").Value = 16950.375
$ws.Range("M135").Value = -13464.2613
$ws.Range("N135").Value = -22020.375
$ws.Range("H141").Value = 1474.5927
$ws.Range("I141").Value = 1223.3636
$ws.Range("J141").Value = 2580
$ws.Range("K141").Value = 3670.0908
$ws.Range("L141").Value = 7740
$ws.Range("N141").Value = -18100
$ws.Range("M141").Value = 1509.9092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3838.79
$ws.Range("I32").Value = 3631.927
$ws.Range("J32").Value = 8803.5
$ws.Range("K32").Value = 3631.927
$ws.Range("L32").Value = 8803.5
$ws.Range("M32").Value = -3344.927
$ws.Range("N32").Value = -9377.5
$ws.Range("H61").Value = 3170.4783
$ws.Range("I61").Value = 3324.8096
$ws.Range("J61").Value = 1550
$ws.Range("K61").Value = 3324.8096
$ws.Range("L61").Value = 1550
$ws.Range("M61").Value = -3112.8096
$ws.Range("N61").Value = -1974
$ws.Range("H110").Value = 486.58334
$ws.Range("I110").Value = 473.125
$ws.Range("J110").Value = 594.25
$ws.Range("K110").Value = 473.125
$ws.Range("L110").Value = 594.25
$ws.Range("M110").Value = 1571.875
$ws.Range("N110").Value = -4684.25
$ws.Range("H132").Value = 5320841.5
$ws.Range("I132").Value = 7354217
$ws.Range("J132").Value = 2781.6924
$ws.Range("K132").Value = 22062651
$ws.Range("L132").Value = 8345.0772
$ws.Range("M132").Value = -22060121
$ws.Range("N132").Value = -13405.0772
$ws.Range("H136").Value = 3170.4783
$ws.Range("I136").Value = 3324.8096
$ws.Range("J136").Value = 1550
$ws.Range("K136").Value = 9974.4288
$ws.Range("L136").Value = 4650
$ws.Range("M136").Value = -7424.4288
$ws.Range("N136").Value = -9750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2341.04
$ws.Range("I20").Value = 2338.3157
$ws.Range("J20").Value = 2349.6667
$ws.Range("K20").Value = 2338.3157
$ws.Range("L20").Value = 2349.6667
$ws.Range("M20").Value = -2091.3157
$ws.Range("N20").Value = -2843.6667
$ws.Range("H86").Value = 45457850
$ws.Range("I86").Value = 2913.25
$ws.Range("J86").Value = 166671000
$ws.Range("K86").Value = 2913.25
$ws.Range("L86").Value = 166671000
$ws.Range("M86").Value = -1790.25
$ws.Range("N86").Value = -166673246
$ws.Range("H89").Value = 45457850
$ws.Range("I89").Value = 2913.25
$ws.Range("J89").Value = 166671000
$ws.Range("K89").Value = 14566.25
$ws.Range("L89").Value = 833355000
$ws.Range("M89").Value = -8950.25
$ws.Range("N89").Value = -833366232
$ws.Range("H99").Value = 874.8570999999999
$ws.Range("I99").Value = 825.4286
$ws.Range("K99").Value = 825.4286
$ws.Range("M99").Value = 672.5714
$ws.Range("H134").Value = 4734.7
$ws.Range("I134").Value = 4114.857
$ws.Range("J134").Value = 5068.4614
$ws.Range("K134").Value = 12344.571
$ws.Range("L134").Value = 15205.3842
$ws.Range("M134").Value = -9809.571
$ws.Range("N134").Value = -20275.3842
$ws.Range("H137").Value = 54449.75
$ws.Range("J137").Value = 54449.75
$ws.Range("L137").Value = 54449.75
$ws.Range("N137").Value = -64649.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 20111
$ws.Range("I104").Value = 16259
$ws.Range("J104").Value = 21395
$ws.Range("K104").Value = 16259
$ws.Range("L104").Value = 21395
$ws.Range("N104").Value = -26637
$ws.Range("M104").Value = -13638
$ws.Range("H122").Value = 1456.4286
$ws.Range("I122").Value = 1531.3334
$ws.Range("J122").Value = 1007
$ws.Range("K122").Value = 4594.0002
$ws.Range("L122").Value = 3021
$ws.Range("M122").Value = -2144.0002
$ws.Range("N122").Value = -7921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 52632420
$ws.Range("I113").Value = 333333950
$ws.Range("J113").Value = 882.375
$ws.Range("K113").Value = 1000001850
$ws.Range("L113").Value = 2647.125
$ws.Range("M113").Value = -999999680
$ws.Range("N113").Value = -6987.125
$ws.Range("H121").Value = 506
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H102").Value = 3711.7827
$ws.Range("I102").Value = 4214.294
$ws.Range("J102").Value = 2288
$ws.Range("K102").Value = 4214.294
$ws.Range("L102").Value = 2288
$ws.Range("M102").Value = -2592.294
$ws.Range("N102").Value = -5532
$ws.Range("H132").Value = 5222.528
$ws.Range("I132").Value = 5627.552
$ws.Range("J132").Value = 3544.5715
$ws.Range("K132").Value = 16882.656
$ws.Range("L132").Value = 10633.7145
$ws.Range("M132").Value = -14352.656
$ws.Range("N132").Value = -15693.7145
$ws.Range("H134").Value = 38800.57
$ws.Range("J134").Value = 38800.57
$ws.Range("L134").Value = 116401.71
$ws.Range("N134").Value = -121471.71
$ws.Range("H136").Value = 33993.332
$ws.Range("J136").Value = 33993.332
$ws.Range("L136").Value = 101979.996
$ws.Range("N136").Value = -107079.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7564.706
$ws.Range("I7").Value = 14750
$ws.Range("K7").Value = 14750
$ws.Range("M7").Value = -14638
$ws.Range("H122").Value = 7289.9546
$ws.Range("I122").Value = 5836.4614
$ws.Range("J122").Value = 9389.444
$ws.Range("K122").Value = 17509.3842
$ws.Range("L122").Value = 28168.332
$ws.Range("M122").Value = -15059.3842
$ws.Range("N122").Value = -33068.33199999999
$ws.Range("H126").Value = 7564.706
$ws.Range("I126").Value = 14750
$ws.Range("K126").Value = 44250
$ws.Range("M126").Value = -41780
$ws.Range("H135").Value = 48112.516
$ws.Range("J135").Value = 48112.516
$ws.Range("L135").Value = 48112.516
$ws.Range("N135").Value = -58252.516

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5187.8
$ws.Range("I100").Value = 10466.5
$ws.Range("J100").Value = 1668.6666
$ws.Range("K100").Value = 20933
$ws.Range("L100").Value = 3337.3332
$ws.Range("M100").Value = -20392
$ws.Range("N100").Value = -4419.3332
$ws.Range("H107").Value = 1941.1428
$ws.Range("I107").Value = 2337.6
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 7012.799999999999
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = -5092.799999999999
$ws.Range("N107").Value = -6690
